$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.426.64"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "3.529.88"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'580.39"
$ws.Range("E5").Value = "  +6.01%  "
$ws.Range("D6").Value = "'179.79"
$ws.Range("E6").Value = "  -5.63%  "
$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  +4.64%  "
$ws.Range("D9").Value = "'0.640"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  +8.74%  "
$ws.Range("D11").Value = "'55.88"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "'0.0000284"
$ws.Range("E12").Value = "  +6.36%  "
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "4.092.40"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").Value = "3.529.54"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'18.48"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.121"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "66.415.69"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").Value = "'12.06"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").Value = "'416.72"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("E22").Value = "  +9.10%  "
$ws.Range("E23").Value = "  +4.05%  "
$ws.Range("D24").Value = "'85.81"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").Value = "'13.24"
$ws.Range("E25").Value = "  +11.41%  "
$ws.Range("D26").Value = "'11.32"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'6.06"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("D29").Value = "'9.16"
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "'6.62"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'604.76"
$ws.Range("E32").Value = "  -7.41%  "
$ws.Range("D33").Value = "'11.75"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "'0.112"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'59.91"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.156"
$ws.Range("E36").Value = "  +11.52%  "
$ws.Range("D37").Value = "0.0₃0814"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'37.40"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "'3.57"
$ws.Range("E40").Value = "  +8.57%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "3.257.11"
$ws.Range("E42").Value = "  +9.05%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "'2.93"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("D45").Value = "'2.57"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.33"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0423"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "'8.69"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "'138.50"
